$d = $word.ActiveDocument
$full = $d.WordOpenXML
$bodyStartTag = "<w:body>"
$bodyStart = $full.IndexOf($bodyStartTag)
$bodyEnd = $full.IndexOf("</w:body>")
$innerStart = $bodyStart + $bodyStartTag.Length
$bodyInner = $full.Substring($innerStart, $bodyEnd - $innerStart)

$d.Content.InsertXML($bodyInner)
Write-Host "After paragraphs count:" $d.Paragraphs.Count

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
Write-Host "Last para text: [" $lastPara.Range.Text "]"
Write-Host "Last para text length:" $lastPara.Range.Text.Length
